$wb = $excel.ActiveWorkbook

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 12347087
$ws.Range("I100").Value = 14493580
$ws.Range("K100").Value = 14493580
$ws.Range("M100").Value = -14493039

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1268
$ws.Range("J125").Value = 1268
$ws.Range("L125").Value = 11412
$ws.Range("N125").Value = -16332

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 16412192
$ws.Range("I132").Value = 3207144
$ws.Range("K132").Value = 9621432
$ws.Range("M132").Value = -9618902

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1042.2632
$ws.Range("I137").Value = 1013.08
$ws.Range("J137").Value = 1098.3846
$ws.Range("K137").Value = 3039.24
$ws.Range("L137").Value = 3295.1538
$ws.Range("M137").Value = -489.2400000000002
$ws.Range("N137").Value = -8395.1538

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 595.8333
$ws.Range("I5").Value = 572.2222
$ws.Range("J5").Value = 666.6667
$ws.Range("K5").Value = 572.2222
$ws.Range("L5").Value = 666.6667
$ws.Range("M5").Value = -460.2222
$ws.Range("N5").Value = -890.6667

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 446796.4
$ws.Range("I32").Value = 4115.0376
$ws.Range("J32").Value = 2401972.5
$ws.Range("K32").Value = 4115.0376
$ws.Range("L32").Value = 2401972.5
$ws.Range("M32").Value = -3828.0376
$ws.Range("N32").Value = -2402546.5

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 595.8333
$ws.Range("I4").Value = 572.2222
$ws.Range("J4").Value = 666.6667
$ws.Range("K4").Value = 572.2222
$ws.Range("L4").Value = 666.6667
$ws.Range("M4").Value = -457.2222
$ws.Range("N4").Value = -896.6667

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1693.8125
$ws.Range("I16").Value = 1573.4
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 1573.4
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -1286.4
$ws.Range("N16").Value = -4074

# CRP row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1782.7273
$ws.Range("I39").Value = 767.7778
$ws.Range("J39").Value = 6350
$ws.Range("K39").Value = 767.7778
$ws.Range("L39").Value = 6350
$ws.Range("M39").Value = -376.7778
$ws.Range("N39").Value = -7132

# CRP row 47
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 7266.6665
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 7266.6665
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 7266.6665
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -8398.666499999999

# CRP row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 1782.7273
$ws.Range("I49").Value = 767.7778
$ws.Range("J49").Value = 6350
$ws.Range("K49").Value = 767.7778
$ws.Range("L49").Value = 6350
$ws.Range("M49").Value = -585.7778
$ws.Range("N49").Value = -6714

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1693.8125
$ws.Range("I113").Value = 1573.4
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1573.4
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 596.5999999999999
$ws.Range("N113").Value = -7840

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7694429.5
$ws.Range("I132").Value = 1959.8572
$ws.Range("J132").Value = 40002804
$ws.Range("K132").Value = 5879.571599999999
$ws.Range("L132").Value = 120008412
$ws.Range("M132").Value = -3349.571599999999
$ws.Range("N132").Value = -120013472

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2000
$ws.Range("J24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("N24").Value = -2346

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1999.9375
$ws.Range("I122").Value = 1999.9333
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5999.7999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3549.7999
$ws.Range("N122").Value = -10900

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2240.818
$ws.Range("I40").Value = 2073.2
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 2073.2
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -1937.2
$ws.Range("N40").Value = -2872

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2356.6667
$ws.Range("I61").Value = 1325
$ws.Range("J61").Value = 2991.5386
$ws.Range("K61").Value = 1325
$ws.Range("L61").Value = 2991.5386
$ws.Range("M61").Value = -1123
$ws.Range("N61").Value = -3395.5386

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2600.8333
$ws.Range("I68").Value = 2033.3334
$ws.Range("J68").Value = 3168.3333
$ws.Range("K68").Value = 2033.3334
$ws.Range("L68").Value = 3168.3333
$ws.Range("M68").Value = -1284.3334
$ws.Range("N68").Value = -4666.3333

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2600.8333
$ws.Range("I71").Value = 2033.3334
$ws.Range("J71").Value = 3168.3333
$ws.Range("K71").Value = 10166.667
$ws.Range("L71").Value = 15841.6665
$ws.Range("M71").Value = -6422.666999999999
$ws.Range("N71").Value = -23329.6665

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2356.6667
$ws.Range("I113").Value = 1325
$ws.Range("J113").Value = 2991.5386
$ws.Range("K113").Value = 1325
$ws.Range("L113").Value = 2991.5386
$ws.Range("M113").Value = 845
$ws.Range("N113").Value = -7331.5386

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3510433.8
$ws.Range("I132").Value = 6667557
$ws.Range("J132").Value = 2519.2222
$ws.Range("K132").Value = 20002671
$ws.Range("L132").Value = 7557.6666
$ws.Range("M132").Value = -20000141
$ws.Range("N132").Value = -12617.6666

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 58469.723
$ws.Range("I136").Value = 83787.914
$ws.Range("K136").Value = 251363.742
$ws.Range("M136").Value = -248813.742

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5971.4062
$ws.Range("I81").Value = 2483.6428
$ws.Range("J81").Value = 8684.111000000001
$ws.Range("K81").Value = 4967.2856
$ws.Range("L81").Value = 17368.222
$ws.Range("M81").Value = -3906.2856
$ws.Range("N81").Value = -19490.222

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5971.4062
$ws.Range("I84").Value = 2483.6428
$ws.Range("J84").Value = 8684.111000000001
$ws.Range("K84").Value = 24836.428
$ws.Range("L84").Value = 86841.11000000002
$ws.Range("M84").Value = -19532.428
$ws.Range("N84").Value = -97449.11000000002

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5270.4814
$ws.Range("I96").Value = 5641.0454
$ws.Range("J96").Value = 3640
$ws.Range("K96").Value = 5641.0454
$ws.Range("L96").Value = 3640
$ws.Range("M96").Value = -4268.0454
$ws.Range("N96").Value = -6386
